# Updated cryptos list data (price + 1h volume change columns, plus a
# 3-row reorder of SuiNetwork/FirstDigitalUSD/ImmutableX in B35:E37).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells get purely-numeric-looking replacement text (e.g. "0.999").
# Pre-mark them as Text so Excel keeps the value as a string instead of
# silently converting it to a number (matches the original inline-string cells).
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '56.873.22'
$ws.Range('E2').Value = '  +4.93%  '
$ws.Range('D3').Value = '2.340.67'
$ws.Range('E3').Value = '  +3.42%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '517.93'
$ws.Range('E5').Value = '  +4.52%  '
$ws.Range('D6').Value = '134.09'
$ws.Range('E6').Value = '  +4.16%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '0.536'
$ws.Range('E8').Value = '  +2.22%  '
$ws.Range('D9').Value = '2.337.86'
$ws.Range('E9').Value = '  +3.06%  '
$ws.Range('D10').Value = '0.103'
$ws.Range('E10').Value = '  +8.72%  '
$ws.Range('E11').Value = '  +1.09%  '
$ws.Range('E12').Value = '  +6.91%  '
$ws.Range('D13').Value = '0.343'
$ws.Range('E13').Value = '  +2.66%  '
$ws.Range('D14').Value = '23.99'
$ws.Range('E14').Value = '  +4.69%  '
$ws.Range('D15').Value = '2.735.06'
$ws.Range('E15').Value = '  +2.67%  '
$ws.Range('D16').Value = '56.765.23'
$ws.Range('E16').Value = '  +4.74%  '
$ws.Range('D17').Value = '0.0000135'
$ws.Range('E17').Value = '  +4.74%  '
$ws.Range('D18').Value = '2.337.62'
$ws.Range('E18').Value = '  +3.07%  '
$ws.Range('D19').Value = '10.54'
$ws.Range('E19').Value = '  +3.26%  '
$ws.Range('E20').Value = '  +4.10%  '
$ws.Range('D21').Value = '321.28'
$ws.Range('E21').Value = '  +6.09%  '
$ws.Range('D22').Value = '6.64'
$ws.Range('E22').Value = '  +5.09%  '
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('D24').Value = '61.24'
$ws.Range('E24').Value = '  +1.04%  '
$ws.Range('D25').Value = '0.995'
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('E26').Value = '  +6.63%  '
$ws.Range('D27').Value = '7.68'
$ws.Range('E27').Value = '  +5.80%  '
$ws.Range('D28').Value = '171.72'
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  +12.22%  '
$ws.Range('E30').Value = '  +6.94%  '
$ws.Range('E31').Value = '  +5.33%  '
$ws.Range('E32').Value = '  +4.72%  '
$ws.Range('D33').Value = '18.38'
$ws.Range('E33').Value = '  +3.48%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.21%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '1.27'
$ws.Range('E36').Value = '  +5.92%  '
$ws.Range('B37').Value = 'SuiNetwork'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D37').Value = '0.943'
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('D38').Value = '4.03'
$ws.Range('E38').Value = '  +8.94%  '
$ws.Range('E39').Value = '  +9.07%  '
$ws.Range('E40').Value = '  +4.70%  '
$ws.Range('D41').Value = '0.382'
$ws.Range('E41').Value = '  +2.05%  '
$ws.Range('D42').Value = '140.79'
$ws.Range('E42').Value = '  +13.07%  '
$ws.Range('E43').Value = '  +7.30%  '
$ws.Range('D44').Value = '5.19'
$ws.Range('E44').Value = '  +8.16%  '
$ws.Range('D45').Value = '276.23'
$ws.Range('E45').Value = '  +14.77%  '
$ws.Range('D46').Value = '0.0510'
$ws.Range('E46').Value = '  +3.70%  '
$ws.Range('E47').Value = '  +4.06%  '
$ws.Range('D48').Value = '0.561'
$ws.Range('E48').Value = '  +3.13%  '
$ws.Range('D49').Value = '0.0216'
$ws.Range('E49').Value = '  +6.11%  '
$ws.Range('E50').Value = '  +2.20%  '
$ws.Range('D51').Value = '16.93'
$ws.Range('E51').Value = '  +5.30%  '
